# Apply data refresh to the two sheets that hold the duplicated event
# listings ("展览" and "全部类型"): bump a couple of "想去人数" /
# "最低票价" figures to reflect newly generated output.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G3").Value = 60
    $ws.Range("F6").Value = 431
    $ws.Range("F9").Value = 541
}
